$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attendance")

# Fix cell B2 (shared string "xc" -> "x c")
$ws.Range("B2").Value = "x c"

# Add new member rows
$ws.Range("A4").Value = 135693201
$ws.Range("B4").Value = "monty python"
$ws.Range("H4").Value = "12:35 PM"

$ws.Range("A5").Value = 789456321
$ws.Range("B5").Value = "safe way"
$ws.Range("H5").Value = "12:36 PM"

$ws.Range("A6").Value = 102450690
$ws.Range("B6").Value = "taylormans"
$ws.Range("H6").Value = "12:48 PM"

$ws.Range("B2").Select()
